$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains literal text formatting (values like "1.000"
# or "0.9999" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.915.43"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.649.54"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "310.85"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "0.3890"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "0.3837"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "51.17"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "0.08437"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "23.82"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "7.013"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "7.963"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "0.00001313"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "1.649.84"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "94.00"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "0.06975"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "19.52"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "6.949"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "13.62"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").Value = "23.916.93"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "2.920"
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("D27").Value = "21.95"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "154.14"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "5.387"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "137.15"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "7.718"
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").Value = "2.486"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "1.830.87"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "0.08162"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "0.9920"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("D36").Value = "6.710"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "0.02914"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("D38").Value = "0.2678"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").Value = "10.47"
$ws.Range("E39").Value = "  -5.65%  "
$ws.Range("D40").Value = "0.09110"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "0.7554"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").Value = "13.47"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "1.423"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "16.84"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").Value = "0.6927"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "2.444"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").Value = "4.098"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("D49").Value = "0.08263"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "134.38"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "1.221"
$ws.Range("E51").Value = "  -2.19%  "
